# Update odds data on Sheet1 to reflect the latest Betfair Back/Lay snapshot
# for 2025-11-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (US MLS - San Diego FC vs Minnesota Utd)
$ws.Range("F2").Value = 1.68
$ws.Range("G2").Value = 1.72
$ws.Range("I2").Value = 5.5
$ws.Range("L2").Value = 1.3
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 5.1
$ws.Range("P2").Value = 2.44
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 2.62
$ws.Range("T2").Value = 1.71
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 1.22
$ws.Range("W2").Value = 2.38
$ws.Range("Z2").Value = 46
$ws.Range("AC2").Value = 10.5
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 18
$ws.Range("AN2").Value = 7.6

# Row 3 (Egyptian Premier - Pyramids vs Al Mokawloon)
$ws.Range("F3").Value = 1.4
$ws.Range("H3").Value = 9.4
$ws.Range("I3").Value = 11.5
$ws.Range("L3").Value = 1.41
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.4
$ws.Range("P3").Value = 1.82
$ws.Range("Q3").Value = 2.02
$ws.Range("R3").Value = 1.31
$ws.Range("S3").Value = 3.55
$ws.Range("T3").Value = 2.24
$ws.Range("U3").Value = 1.65
$ws.Range("X3").Value = 17.5
$ws.Range("Y3").Value = 980
$ws.Range("Z3").Value = 120
$ws.Range("AA3").Value = 570
$ws.Range("AB3").Value = 7
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 980
$ws.Range("AE3").Value = 260
$ws.Range("AF3").Value = 7.8
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 220
$ws.Range("AJ3").Value = 12.5
$ws.Range("AK3").Value = 18.5
$ws.Range("AL3").Value = 55
$ws.Range("AM3").Value = 290
$ws.Range("AN3").Value = 8.800000000000001
$ws.Range("AO3").Value = 470

# Row 4 (Dutch Eerste Divisie - Vitesse Arnhem vs Jong PSV Eindhoven)
$ws.Range("P4").Value = 2.72

# Row 6 (Dutch Eerste Divisie - Roda JC vs FC Dordrecht)
$ws.Range("J6").Value = 3.6

# Row 7 (Dutch Eerste Divisie - RKC Waalwijk vs MVV Maastricht)
$ws.Range("J7").Value = 4.5
